$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 21.797493
$ws.Range("N2").Value = 65.39247899999999
$ws.Range("O2").Value = 0.5076406019363033
$ws.Range("P2").Value = 0.5076406019363033
$ws.Range("Q2").Value = 5.080152781904
$ws.Range("R2").Value = 45.721375037136
$ws.Range("S2").Value = 0.3611572878511963
$ws.Range("T2").Value = 0.3611572878511963
$ws.Range("O3").Value = 0.1566265815027233
$ws.Range("P3").Value = 0.1566265815027233
$ws.Range("S3").Value = 0.1114308649961487
$ws.Range("T3").Value = 0.1114308649961487
$ws.Range("M4").Value = 6.908085666666667
$ws.Range("N4").Value = 20.724257
$ws.Range("O4").Value = 0.1608820228112571
$ws.Range("P4").Value = 0.1608820228112571
$ws.Range("Q4").Value = 1.610007656254222
$ws.Range("R4").Value = 14.490068906288
$ws.Range("S4").Value = 0.1144583683828712
$ws.Range("T4").Value = 0.1144583683828712
$ws.Range("M5").Value = 7.507888333333334
$ws.Range("N5").Value = 22.523665
$ws.Range("O5").Value = 0.1748507937497162
$ws.Range("P5").Value = 0.1748507937497162
$ws.Range("Q5").Value = 1.749798465484445
$ws.Range("R5").Value = 15.74818618936
$ws.Range("S5").Value = 0.12439635089945
$ws.Range("T5").Value = 0.12439635089945
$ws.Range("M6").Value = 21.797493
$ws.Range("N6").Value = 65.39247899999999
$ws.Range("O6").Value = 0.5076406019363033
$ws.Range("P6").Value = 0.5076406019363033
$ws.Range("Q6").Value = 2.060480684135
$ws.Range("R6").Value = 18.544326157215
$ws.Range("S6").Value = 0.1464833140851071
$ws.Range("T6").Value = 0.1464833140851071
$ws.Range("O7").Value = 0.1566265815027233
$ws.Range("P7").Value = 0.1566265815027233
$ws.Range("S7").Value = 0.04519571650657458
$ws.Range("T7").Value = 0.04519571650657458
$ws.Range("M8").Value = 6.908085666666667
$ws.Range("N8").Value = 20.724257
$ws.Range("O8").Value = 0.1608820228112571
$ws.Range("P8").Value = 0.1608820228112571
$ws.Range("Q8").Value = 0.6530098245938889
$ws.Range("R8").Value = 5.877088421345
$ws.Range("S8").Value = 0.04642365442838586
$ws.Range("T8").Value = 0.04642365442838586
$ws.Range("M9").Value = 7.507888333333334
$ws.Range("N9").Value = 22.523665
$ws.Range("O9").Value = 0.1748507937497162
$ws.Range("P9").Value = 0.1748507937497162
$ws.Range("Q9").Value = 0.7097081710027777
$ws.Range("R9").Value = 6.387373539025
$ws.Range("S9").Value = 0.05045444285026621
$ws.Range("T9").Value = 0.05045444285026621
